$d = $word.ActiveDocument

# 1. Main body paragraph text changes
$d.Content.Find.Execute(
    "USDT 帳戶將於 2023 年 9 月 29 日格林尼治標準時間 00:00 關閉。 任何未平倉頭寸將自動關閉，帳戶餘額將在上述日期後轉移到您最後一個有效的帳戶中",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "USDT 帳戶將於 2023 年 9 月 29 日格林尼治標準時間 00:00 關閉。 任何持倉頭寸將在上述日期後自動平倉，帳戶餘額將轉移到最後活躍的帳戶",
    2
)

$d.Content.Find.Execute(
    "在此過程中將適用標準匯率和費用。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "在此過程中將採用標準匯率和費用。",
    2
)

$d.Content.Find.Execute(
    "如有任何疑問，請聯繫我們：",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "如有任何疑問，請透過以下方式聯繫我們：",
    2
)

# 2. Comment text changes
$d.Comments(1).Range.Find.Execute(
    "@azita@regentmarkets.com，BE 不能保證他們能在那時準備好腳本，",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "@azita@regentmarkets.com，BE 無法保證屆時可以完成指令，",
    2
)

$d.Comments(1).Range.Find.Execute(
    "我們可以說成「在提到的日期之後」嗎？",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "可以將其稱為“所述日期之後”嗎？",
    2
)

$d.Comments(2).Range.Find.Execute(
    "您是指轉帳將在提到的日期之後進行嗎？",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "意思是在所述日期之後完成轉帳嗎？",
    2
)

$d.Comments(3).Range.Find.Execute(
    "是的..看起來我們無法確認日期",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "是的... 看來我們無法確認日期",
    2
)
